# Update the cryptos list with refreshed prices / 1h volume percentages.
# Values such as "1.012" look numeric to Excel, so they are entered with a
# leading apostrophe to force text storage and avoid losing significant
# trailing zeros (matches the sheet's original inline-string / text layout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.361.99"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "1.833.35"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "'314.04"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "'0.4743"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.07458"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'0.8860"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").Value = "'20.46"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.899.33"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").Value = "'0.07344"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("D14").Value = "'5.453"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'93.10"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "'6.586"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'0.000008824"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "27.679.67"
$ws.Range("D21").Value = "'14.79"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'5.312"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'10.68"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "2.102.53"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").Value = "'1.890"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "'151.94"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "'2.140"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "'5.242"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'117.49"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").Value = "'0.09004"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").Value = "'0.7549"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "'4.550"
$ws.Range("D35").Value = "'2.941"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'1.011"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "'1.103"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "'0.05354"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D40").Value = "'2.983"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "'7.321"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'2.401"
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("D43").Value = "'0.5324"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'0.1661"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").Value = "'8.497"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'0.4917"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "'10.58"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.10"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.011"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'0.06301"
$ws.Range("E51").Value = "  +0.06%  "
